$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 400
$ws.Range("B2").Value = 387
$ws.Range("C2").Value = 13

$ws.Range("B5").Value = 0.9675
$ws.Range("C5").Value = 0.0325
